$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New week of price data was inserted at the top (rows 2-3). Every existing
# record shifts down by two rows; the two oldest records that fall off the
# bottom of the table are preserved by appending two brand-new rows at the
# end (110-111), growing the sheet's used range from A1:T109 to A1:T111.
# Only columns D (Fecha), L (Calidad), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), Q (Unidad de
# comercializacion), R (Origen), S (Precio $/Kg) and T (Kg/unidad) actually
# vary row to row - A, B, C, E, F, G, H, I, J, K are constant for every
# record in this subset, so they are left untouched for rows 2-109.
# ---------------------------------------------------------------------------

$dataRows = @'
2|44956|Primera|100|3000|3000|3000|$/bandeja 2 kilos|Provincia de Curicó|1500|2
3|44956|Segunda|100|2500|2500|2500|$/bandeja 2 kilos|Provincia de Curicó|1250|2
4|44585|Segunda|150|3000|3000|3000|$/bandeja 2 kilos|Provincia de Linares|1500|2
5|44925|Primera|150|3000|3000|3000|$/bandeja 2 kilos|Provincia de Curicó|1500|2
6|44567|Primera|180|3600|3600|3600|$/bandeja 2 kilos|Provincia de Linares|1800|2
7|44172|Primera|300|3400|3600|3467|$/bandeja 2 kilos|Provincia de Linares|1734|2
8|44937|Primera|20|3000|3000|3000|$/bandeja 2 kilos|Provincia de Curicó|1500|2
9|44910|Primera|150|3000|3000|3000|$/bandeja 2 kilos|Provincia de Curicó|1500|2
10|44554|Primera|100|3800|3800|3800|$/bandeja 2 kilos|Provincia de Linares|1900|2
11|44624|Primera|120|3300|3300|3300|$/bandeja 2 kilos|Provincia de Linares|1650|2
12|44624|Segunda|100|3000|3000|3000|$/bandeja 2 kilos|Provincia de Linares|1500|2
13|44186|Primera|200|3000|3000|3000|$/bandeja 2 kilos|Provincia de Limarí|1500|2
14|44895|Primera|150|3600|3600|3600|$/bandeja 2 kilos|Provincia de Curicó|1800|2
15|44895|Primera|220|3600|4000|3782|$/bandeja 2 kilos|Provincia de Curicó|1891|2
16|44265|Primera|70|3600|3800|3714|$/bandeja 2 kilos|Provincia de Linares|1857|2
17|44511|Primera|50|6400|6400|6400|$/bandeja 2 kilos|Provincia de Linares|3200|2
18|44930|Primera|170|2800|3000|2859|$/bandeja 2 kilos|Provincia de Curicó|1430|2
19|44537|Primera|200|3600|3600|3600|$/bandeja 2 kilos|Provincia de Linares|1800|2
20|44537|Segunda|100|3000|3000|3000|$/bandeja 2 kilos|Provincia de Linares|1500|2
21|44525|Primera|150|4000|4000|4000|$/bandeja 2 kilos|Provincia de Linares|2000|2
22|44571|Segunda|120|3200|3200|3200|$/bandeja 2 kilos|Provincia de Linares|1600|2
23|44659|Segunda|60|3000|3000|3000|$/bandeja 2 kilos|Provincia de Linares|1500|2
24|44532|Primera|170|3600|3600|3600|$/bandeja 2 kilos|Provincia de Linares|1800|2
25|44918|Primera|100|2800|3000|2900|$/bandeja 2 kilos|Provincia de Colchagua|1450|2
26|44187|Primera|110|2600|3000|2782|$/bandeja 2 kilos|Provincia de Linares|1391|2
27|44210|Segunda|150|2700|2700|2700|$/bandeja 2 kilos|Provincia de Linares|1350|2
28|44544|Primera|200|4000|4000|4000|$/bandeja 2 kilos|Provincia de Linares|2000|2
29|44547|Primera|150|3600|3600|3600|$/bandeja 2 kilos|Provincia de Linares|1800|2
30|44547|Segunda|100|3000|3000|3000|$/bandeja 2 kilos|Provincia de Linares|1500|2
31|44522|Primera|30|5000|5000|5000|$/bandeja 2 kilos|Provincia de Linares|2500|2
32|44876|Primera|200|6000|6000|6000|$/bandeja 2 kilos|Provincia de Curicó|3000|2
33|44936|Primera|40|3000|3000|3000|$/bandeja 2 kilos|Provincia de Curicó|1500|2
34|44893|Primera|100|4000|4000|4000|$/bandeja 2 kilos|Provincia de Curicó|2000|2
35|44516|Primera|30|5000|5000|5000|$/bandeja 2 kilos|Provincia de Linares|2500|2
36|44655|Segunda|60|3000|3000|3000|$/bandeja 2 kilos|Provincia de Linares|1500|2
37|44622|Primera|80|3000|3000|3000|$/bandeja 2 kilos|Provincia de Linares|1500|2
38|44596|Primera|200|3500|3500|3500|$/bandeja 2 kilos|Provincia de Linares|1750|2
39|44873|Primera|35|3500|3500|3500|$/bandeja 2 kilos|Provincia de Curicó|1750|2
40|44952|Primera|150|3000|3000|3000|$/bandeja 2 kilos|Provincia de Curicó|1500|2
41|44952|Segunda|200|2600|2600|2600|$/bandeja 2 kilos|Provincia de Curicó|1300|2
42|44907|Primera|180|3000|3000|3000|$/bandeja 2 kilos|Provincia de Curicó|1500|2
43|44890|Primera|100|3600|3600|3600|$/bandeja 2 kilos|Provincia de Curicó|1800|2
44|44573|Primera|120|3600|3600|3600|$/bandeja 2 kilos|Provincia de Linares|1800|2
45|44530|Primera|100|4000|4000|4000|$/bandeja 2 kilos|Provincia de Linares|2000|2
46|44530|Segunda|150|3600|3600|3600|$/bandeja 2 kilos|Provincia de Linares|1800|2
47|44589|Primera|150|3500|3500|3500|$/bandeja 2 kilos|Provincia de Linares|1750|2
48|44162|Primera|100|4000|4000|4000|$/bandeja 2 kilos|Región de O'Higgins|2000|2
49|44533|Primera|180|3600|3600|3600|$/bandeja 2 kilos|Provincia de Linares|1800|2
50|44533|Segunda|100|3000|3000|3000|$/bandeja 2 kilos|Provincia de Linares|1500|2
51|44545|Primera|150|3800|3800|3800|$/bandeja 2 kilos|Provincia de Linares|1900|2
52|44519|Primera|180|4000|4000|4000|$/bandeja 2 kilos|Provincia de Linares|2000|2
53|44512|Primera|30|6000|6000|6000|$/bandeja 2 kilos|Provincia de Linares|3000|2
54|44932|Primera|180|3000|3000|3000|$/bandeja 2 kilos|Provincia de Curicó|1500|2
55|44167|Primera|500|3600|3600|3600|$/bandeja 2 kilos|Región de O'Higgins|1800|2
56|44232|Primera|60|3000|3000|3000|$/bandeja 2 kilos|Provincia de Linares|1500|2
57|44939|Primera|45|3000|3000|3000|$/bandeja 2 kilos|Provincia de Curicó|1500|2
58|44939|Segunda|30|2600|2600|2600|$/bandeja 2 kilos|Provincia de Curicó|1300|2
59|44641|Segunda|50|3000|3000|3000|$/bandeja 2 kilos|Provincia de Linares|1500|2
60|44211|Primera|40|2800|2800|2800|$/bandeja 2 kilos|Provincia de Linares|1400|2
61|44211|Segunda|30|2600|2600|2600|$/bandeja 2 kilos|Provincia de Linares|1300|2
62|44546|Primera|100|3800|3800|3800|$/bandeja 2 kilos|Provincia de Linares|1900|2
63|44546|Segunda|150|3000|3000|3000|$/bandeja 2 kilos|Provincia de Linares|1500|2
64|44917|Primera|60|3000|3000|3000|$/bandeja 2 kilos|Provincia de Curicó|1500|2
65|44517|Primera|20|5000|5000|5000|$/bandeja 2 kilos|Provincia de Linares|2500|2
66|44635|Segunda|120|3000|3000|3000|$/bandeja 2 kilos|Provincia de Linares|1500|2
67|44536|Primera|180|3600|3600|3600|$/bandeja 2 kilos|Provincia de Linares|1800|2
68|44536|Segunda|80|3000|3000|3000|$/bandeja 2 kilos|Provincia de Linares|1500|2
69|44883|Primera|250|5600|6000|5760|$/bandeja 2 kilos|Región de O'Higgins|2880|2
70|44176|Primera|150|3500|3500|3500|$/bandeja 12 canastillos 125 gramos|Provincia de Curicó|2333|1.5
71|44165|Primera|400|3400|3400|3400|$/bandeja 2 kilos|Región de O'Higgins|1700|2
72|44915|Primera|220|3000|3200|3109|$/bandeja 2 kilos|Provincia de Curicó|1554|2
73|44588|Primera|150|3500|3500|3500|$/bandeja 2 kilos|Provincia de Linares|1750|2
74|44202|Primera|30|3000|3000|3000|$/bandeja 2 kilos|Provincia de Linares|1500|2
75|44202|Segunda|20|2600|2600|2600|$/bandeja 2 kilos|Provincia de Linares|1300|2
76|44931|Primera|150|3000|3000|3000|$/bandeja 2 kilos|Provincia de Curicó|1500|2
77|44894|Primera|220|3600|4000|3782|$/bandeja 2 kilos|Provincia de Curicó|1891|2
78|44938|Primera|30|3000|3000|3000|$/bandeja 2 kilos|Provincia de Curicó|1500|2
79|44200|Segunda|50|2600|2600|2600|$/bandeja 2 kilos|Provincia de Linares|1300|2
80|44518|Primera|20|5000|5000|5000|$/bandeja 2 kilos|Provincia de Linares|2500|2
81|44204|Primera|50|3000|3000|3000|$/bandeja 2 kilos|Provincia de Linares|1500|2
82|44204|Segunda|140|2400|2400|2400|$/bandeja 2 kilos|Provincia de Linares|1200|2
83|44902|Primera|250|3000|3200|3080|$/bandeja 2 kilos|Provincia de Curicó|1540|2
84|44526|Primera|250|4000|4000|4000|$/bandeja 2 kilos|Provincia de Linares|2000|2
85|44586|Segunda|150|3000|3000|3000|$/bandeja 2 kilos|Provincia de Linares|1500|2
86|44166|Primera|1500|3600|3600|3600|$/bandeja 2 kilos|Región de O'Higgins|1800|2
87|44582|Segunda|160|3200|3200|3200|$/bandeja 2 kilos|Provincia de Linares|1600|2
88|44900|Primera|270|3000|3000|3000|$/bandeja 2 kilos|Provincia de Curicó|1500|2
89|44524|Primera|180|4000|4000|4000|$/bandeja 2 kilos|Provincia de Linares|2000|2
90|44882|Primera|170|6000|6000|6000|$/bandeja 2 kilos|Región de O'Higgins|3000|2
91|44587|Primera|150|3500|3500|3500|$/bandeja 2 kilos|Provincia de Linares|1750|2
92|44264|Primera|110|3500|4000|3727|$/bandeja 2 kilos|Provincia de Linares|1864|2
93|44951|Segunda|40|2600|2600|2600|$/bandeja 2 kilos|Provincia de Curicó|1300|2
94|44935|Primera|40|3000|3000|3000|$/bandeja 2 kilos|Provincia de Curicó|1500|2
95|44904|Primera|150|3000|3000|3000|$/bandeja 2 kilos|Provincia de Curicó|1500|2
96|44897|Primera|210|3000|3000|3000|$/bandeja 2 kilos|Provincia de Curicó|1500|2
97|44529|Primera|100|3800|3800|3800|$/bandeja 2 kilos|Provincia de Linares|1900|2
98|44911|Primera|220|3000|3000|3000|$/bandeja 2 kilos|Provincia de Curicó|1500|2
99|44901|Primera|230|3000|3000|3000|$/bandeja 2 kilos|Provincia de Curicó|1500|2
100|44946|Primera|50|3000|3000|3000|$/bandeja 2 kilos|Provincia de Curicó|1500|2
101|44637|Segunda|150|3000|3000|3000|$/bandeja 2 kilos|Provincia de Linares|1500|2
102|44581|Segunda|200|3400|3400|3400|$/bandeja 2 kilos|Provincia de Linares|1700|2
103|44631|Segunda|60|3700|3700|3700|$/bandeja 2 kilos|Provincia de Linares|1850|2
104|44908|Primera|250|3000|3000|3000|$/bandeja 2 kilos|Provincia de Curicó|1500|2
105|44235|Primera|60|3000|3000|3000|$/bandeja 2 kilos|Provincia de Linares|1500|2
106|44539|Segunda|150|3000|3000|3000|$/bandeja 2 kilos|Provincia de Linares|1500|2
107|44592|Primera|150|3500|3500|3500|$/bandeja 2 kilos|Provincia de Linares|1750|2
108|44169|Primera|400|3600|3600|3600|$/bandeja 2 kilos|Provincia de Linares|1800|2
109|44552|Primera|180|4000|4000|4000|$/bandeja 2 kilos|Provincia de Linares|2000|2
'@

$newRows = @'
110|5|Macroferia Regional de Talca|Maule|44595|7|Fruta|100101|Berries|100101001|Arándano (blue)|Sin especificar|Primera|150|3500|3500|3500|$/bandeja 2 kilos|Provincia de Linares|1750|2
111|5|Macroferia Regional de Talca|Maule|44889|7|Fruta|100101|Berries|100101001|Arándano (blue)|Sin especificar|Primera|150|4000|4000|4000|$/bandeja 2 kilos|Región de O'Higgins|2000|2
'@

foreach ($line in $dataRows -split "`n") {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $f = $line -split "\|"
    $r = [int]$f[0]

    $ws.Cells.Item($r, 4).Value  = [double]$f[1]     # D - Fecha
    $ws.Cells.Item($r, 12).Value = $f[2]              # L - Calidad
    $ws.Cells.Item($r, 13).Value = [double]$f[3]     # M - Volumen
    $ws.Cells.Item($r, 14).Value = [double]$f[4]     # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = [double]$f[5]     # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = [double]$f[6]     # P - Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $f[7]              # Q - Unidad de comercializacion
    $ws.Cells.Item($r, 18).Value = $f[8]              # R - Origen
    $ws.Cells.Item($r, 19).Value = [double]$f[9]     # S - Precio $/Kg
    $ws.Cells.Item($r, 20).Value = [double]$f[10]    # T - Kg / unidad
}

foreach ($line in $newRows -split "`n") {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $f = $line -split "\|"
    $r = [int]$f[0]

    $ws.Cells.Item($r, 1).Value  = [double]$f[1]     # A - Mercado ID
    $ws.Cells.Item($r, 2).Value  = $f[2]              # B - Mercado
    $ws.Cells.Item($r, 3).Value  = $f[3]              # C - Region
    $ws.Cells.Item($r, 4).Value  = [double]$f[4]     # D - Fecha
    $ws.Cells.Item($r, 5).Value  = [double]$f[5]     # E - Codreg
    $ws.Cells.Item($r, 6).Value  = $f[6]              # F - Tipo
    $ws.Cells.Item($r, 7).Value  = [double]$f[7]     # G - Producto ID
    $ws.Cells.Item($r, 8).Value  = $f[8]              # H - Producto
    $ws.Cells.Item($r, 9).Value  = [double]$f[9]     # I - Categoria ID
    $ws.Cells.Item($r, 10).Value = $f[10]             # J - Categoria
    $ws.Cells.Item($r, 11).Value = $f[11]             # K - Variedad
    $ws.Cells.Item($r, 12).Value = $f[12]             # L - Calidad
    $ws.Cells.Item($r, 13).Value = [double]$f[13]    # M - Volumen
    $ws.Cells.Item($r, 14).Value = [double]$f[14]    # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = [double]$f[15]    # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = [double]$f[16]    # P - Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $f[17]             # Q - Unidad de comercializacion
    $ws.Cells.Item($r, 18).Value = $f[18]             # R - Origen
    $ws.Cells.Item($r, 19).Value = [double]$f[19]    # S - Precio $/Kg
    $ws.Cells.Item($r, 20).Value = [double]$f[20]    # T - Kg / unidad

    # Match the date-style formatting used by the "Fecha" column elsewhere
    # in the table (cell style carries the YYYY-MM-DD HH:MM:SS number format).
    $ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item($r - 1, 4).NumberFormat
}
